$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "H 72" record (originally row 2); all subsequent rows shift up by one.
$ws.Rows.Item(2).Delete()
